# relatorio2.xlsx -- "Trabalho final com parte de SO"
# Expand the 1-file-vs-method comparison into a 10-file comparison table
# (Arquivo 10..20) with one series per scheduling method, add a chart
# title/legend, and move the chart below the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B already carries the "bestFit" width computed for the original
# table; reuse that same width for every new data column (B:L).
$dataColWidth = $ws.Range("B1").ColumnWidth

# ---------------------------------------------------------------------
# 1) Re-shape the data table.
#    Old layout:            New layout:
#      A1 Método | B1 Tempo(ms)     A1 Método(merged A1:A2) | B1..L1 Arquivo 10..20
#      A2 Fila...| B2 164           A2 (merged)             | B2..L2 Tempo (ms)
#      A3 Menor..| B3 78            A3 Fila de Prioridade    | B3..K3 values
#      A4 Ordem..| B4 97            A4 Menor Primeiro        | B4..K4 values
#      A5 Round..| B5 89            A5 Ordem Chegada         | B5..K5 values
#                                    A6 Round Robin           | B6..K6 values
# ---------------------------------------------------------------------

# Insert a new row at the top so the existing method rows shift down
# from 2..5 to 3..6, then build the two header rows.
$ws.Rows("1:1").Insert()

$ws.Range("A1").Value = "Método"

$fileNames = @("Arquivo 10","Arquivo 11","Arquivo 12","Arquivo 13","Arquivo 14","Arquivo 15","Arquivo 16","Arquivo 17","Arquivo 18","Arquivo 19","Arquivo 20")
for ($i = 0; $i -lt $fileNames.Length; $i++) {
    $col = 2 + $i
    $ws.Cells.Item(1, $col).Value = $fileNames[$i]
}

# Row 2: repeat "Tempo (ms)" under every file column (B..L).
for ($i = 0; $i -lt $fileNames.Length; $i++) {
    $col = 2 + $i
    $ws.Cells.Item(2, $col).Value = "Tempo (ms)"
}

# Data rows 3..6 (C..K, columns 3-11) for the four scheduling methods.
$priFila = @(11,9,9,9,24,37,54,54,193)
$menor   = @(10,6,6,11,22,31,39,61,171)
$ordem   = @(0,1,2,1,1,2,2,4,12)
$round   = @(17,31,63,143,578,2106,9028,55591,196396)

$ws.Range("B3").Value = 82
$ws.Range("B4").Value = 4
$ws.Range("B5").Value = 1
$ws.Range("B6").Value = 31

for ($i = 0; $i -lt $priFila.Length; $i++) {
    $col = 3 + $i
    $ws.Cells.Item(3, $col).Value = $priFila[$i]
    $ws.Cells.Item(4, $col).Value = $menor[$i]
    $ws.Cells.Item(5, $col).Value = $ordem[$i]
    $ws.Cells.Item(6, $col).Value = $round[$i]
}

# ---------------------------------------------------------------------
# 2) Formatting: copy the existing header/data styles onto the new cells.
# ---------------------------------------------------------------------

# B1:L1 + A1 get the bold/filled header look (same visual style the old
# A1/B1 header row used).
$ws.Range("A2").Copy() | Out-Null
$ws.Range("B1:L1").PasteSpecial(-4122) | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null

# B2:L2 ("Tempo (ms)") and A3:A6 (method names) reuse the bold header style too.
$ws.Range("B2:L2").PasteSpecial(-4122) | Out-Null
$ws.Range("A3:A6").PasteSpecial(-4122) | Out-Null

# B3:B6 reuse the centered data style that used to live on column B.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B4:B6").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Merge the "Método" header down both rows.
$ws.Range("A1:A2").Merge() | Out-Null
$ws.Range("A1").Value = "Método"

# Column widths: same bestFit width across every data column B:L (mirrors
# column B's original auto-fit width).
$ws.Range("B1:L6").ColumnWidth = $dataColWidth

# ---------------------------------------------------------------------
# 3) Selection / view tweaks to match the saved workbook state.
# ---------------------------------------------------------------------
$ws.Range("B1:L1").Select() | Out-Null

# ---------------------------------------------------------------------
# 4) Chart: retarget the single series onto 4 series (one per method),
#    using the file-name row as categories; add title + legend.
# ---------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

$s1 = $chart.SeriesCollection().Item(1)
$s1.Formula = "=SERIES(Plan1!`$A`$3,Plan1!`$B`$1:`$L`$1,Plan1!`$B`$3:`$L`$3,1)"

$s2 = $chart.SeriesCollection().NewSeries()
$s2.Formula = "=SERIES(Plan1!`$A`$4,Plan1!`$B`$1:`$L`$1,Plan1!`$B`$4:`$L`$4,2)"

$s3 = $chart.SeriesCollection().NewSeries()
$s3.Formula = "=SERIES(Plan1!`$A`$5,Plan1!`$B`$1:`$L`$1,Plan1!`$B`$5:`$L`$5,3)"

$s4 = $chart.SeriesCollection().NewSeries()
$s4.Formula = "=SERIES(Plan1!`$A`$6,Plan1!`$B`$1:`$L`$1,Plan1!`$B`$6:`$L`$6,4)"

$chart.HasTitle = $true
$chart.ChartTitle.Text = "Comparativo"

$chart.HasLegend = $true
$chart.Legend.Position = -4107

# ---------------------------------------------------------------------
# 5) Move/resize the chart to sit below the table (col A, row 8 .. col M, row 22).
# ---------------------------------------------------------------------
$fromCell = $ws.Cells.Item(8, 1)
$toCell = $ws.Cells.Item(22, 13)

$co.Left = $fromCell.Left
$co.Top = $fromCell.Top + 14287 / 12700
$co.Width = ($toCell.Left + 19050 / 12700) - $fromCell.Left
$co.Height = ($toCell.Top + 90487 / 12700) - ($fromCell.Top + 14287 / 12700)

$co.Name = "Gráfico 2"
